$wb = $excel.ActiveWorkbook

# --- Update stale selections left on the two pre-existing sheets that were navigated ---
$ws1 = $wb.Worksheets.Item("Fat pad wts")
$ws1.Range("Q3:Q15").Select()

$ws3 = $wb.Worksheets.Item("Muscle wts")
$ws3.Range("A3:J14").Select()

# --- Add the new "Tissue Weight Summary" sheet (will be moved to the end below) ---
$ws4 = $wb.Worksheets.Add()
$ws4.Name = "Tissue Weight Summary"

# Match the default page margins used by the other sheets in this workbook (values are in points: 1in = 72pt).
$ws4.PageSetup.LeftMargin = 54
$ws4.PageSetup.RightMargin = 54
$ws4.PageSetup.TopMargin = 72
$ws4.PageSetup.BottomMargin = 72
$ws4.PageSetup.HeaderMargin = 36
$ws4.PageSetup.FooterMargin = 36

# --- Prime shared-string insertion order to match target indices ---
$ws4.Range("D1").Value = "Right.EWAT"
$ws4.Range("E1").Value = "Left EWAT"
$ws4.Range("A1").Value = "Diet"
$ws4.Range("B1").Value = "Mouse"
$ws4.Range("A2").Value = "Control Diet"
$ws4.Range("A8").Value = "High Protein Diet"
$ws4.Range("F1").Value = "Right.IWAT"
$ws4.Range("G1").Value = "Left.IWAT"
$ws4.Range("C1").Value = "Body Weight"
$ws4.Range("H1").Value = "Right.Quad"
$ws4.Range("I1").Value = "Left.Quad"
$ws4.Range("J1").Value = "TS.Left"
$ws4.Range("K1").Value = "TS.Right"

# --- Fill in all remaining data (rows 1-13), in natural reading order ---
$ws4.Range("L1").Value = "Heart"
$ws4.Range("B2").Value = 2774
$ws4.Range("C2").Value = 35
$ws4.Range("D2").Value = 1081.9000000000001
$ws4.Range("E2").Value = 1075.2
$ws4.Range("F2").Value = 607.4
$ws4.Range("G2").Value = 676.6
$ws4.Range("H2").Value = 236.9
$ws4.Range("H2").Font.Color = 0
$ws4.Range("I2").Value = 235.4
$ws4.Range("J2").Value = 173.5
$ws4.Range("K2").Value = 187.4
$ws4.Range("L2").Value = 126.5
$ws4.Range("A3").Value = "Control Diet"
$ws4.Range("B3").Value = 2777
$ws4.Range("C3").Value = 40.700000000000003
$ws4.Range("D3").Value = 458.9
$ws4.Range("E3").Value = 548.4
$ws4.Range("F3").Value = 288
$ws4.Range("G3").Value = 270.2
$ws4.Range("H3").Value = 279.3
$ws4.Range("I3").Value = 248.9
$ws4.Range("J3").Value = 182.9
$ws4.Range("K3").Value = 181.5
$ws4.Range("L3").Value = 127.7
$ws4.Range("A4").Value = "Control Diet"
$ws4.Range("B4").Value = 2778
$ws4.Range("C4").Value = 39.799999999999997
$ws4.Range("A5").Value = "Control Diet"
$ws4.Range("B5").Value = 2792
$ws4.Range("C5").Value = 41.2
$ws4.Range("D5").Value = 1214.5999999999999
$ws4.Range("E5").Value = 936.9
$ws4.Range("F5").Value = 603.6
$ws4.Range("G5").Value = 577.9
$ws4.Range("H5").Value = 224
$ws4.Range("I5").Value = 228.4
$ws4.Range("J5").Value = 164.7
$ws4.Range("K5").Value = 133
$ws4.Range("L5").Value = 122.4
$ws4.Range("A6").Value = "Control Diet"
$ws4.Range("B6").Value = 2791
$ws4.Range("C6").Value = 37.5
$ws4.Range("A7").Value = "Control Diet"
$ws4.Range("B7").Value = 2797
$ws4.Range("C7").Value = 30.8
$ws4.Range("D7").Value = 267.7
$ws4.Range("E7").Value = 248.6
$ws4.Range("F7").Value = 185.6
$ws4.Range("G7").Value = 184.5
$ws4.Range("H7").Value = 223.9
$ws4.Range("I7").Value = 241.8
$ws4.Range("J7").Value = 171.9
$ws4.Range("K7").Value = 173.6
$ws4.Range("L7").Value = 129.1
$ws4.Range("B8").Value = 2776
$ws4.Range("C8").Value = 30.1
$ws4.Range("D8").Value = 197.3
$ws4.Range("E8").Value = 213.3
$ws4.Range("F8").Value = 103.2
$ws4.Range("G8").Value = 118.5
$ws4.Range("H8").Value = 258.2
$ws4.Range("I8").Value = 255.9
$ws4.Range("J8").Value = 188.1
$ws4.Range("K8").Value = 173.7
$ws4.Range("L8").Value = 142.1
$ws4.Range("A9").Value = "High Protein Diet"
$ws4.Range("B9").Value = 2784
$ws4.Range("C9").Value = 37.1
$ws4.Range("A10").Value = "High Protein Diet"
$ws4.Range("B10").Value = 2781
$ws4.Range("C10").Value = 31.4
$ws4.Range("D10").Value = 575.4
$ws4.Range("E10").Value = 548.79999999999995
$ws4.Range("F10").Value = 348.8
$ws4.Range("G10").Value = 327.5
$ws4.Range("H10").Value = 265.3
$ws4.Range("I10").Value = 249.4
$ws4.Range("J10").Value = 200.1
$ws4.Range("K10").Value = 195.2
$ws4.Range("L10").Value = 133.6
$ws4.Range("A11").Value = "High Protein Diet"
$ws4.Range("B11").Value = 2790
$ws4.Range("C11").Value = 40.200000000000003
$ws4.Range("D11").Value = 597.4
$ws4.Range("E11").Value = 664.5
$ws4.Range("F11").Value = 401.9
$ws4.Range("G11").Value = 391.5
$ws4.Range("H11").Value = 273.8
$ws4.Range("I11").Value = 249.8
$ws4.Range("J11").Value = 202.2
$ws4.Range("K11").Value = 196.2
$ws4.Range("L11").Value = 141.19999999999999
$ws4.Range("A12").Value = "High Protein Diet"
$ws4.Range("B12").Value = 2795
$ws4.Range("C12").Value = 32.799999999999997
$ws4.Range("D12").Value = 221.7
$ws4.Range("E12").Value = 236.6
$ws4.Range("F12").Value = 168.5
$ws4.Range("G12").Value = 149.69999999999999
$ws4.Range("H12").Value = 262.3
$ws4.Range("I12").Value = 249.6
$ws4.Range("J12").Value = 176.7
$ws4.Range("K12").Value = 181.6
$ws4.Range("L12").Value = 124.6
$ws4.Range("A13").Value = "High Protein Diet"
$ws4.Range("B13").Value = 2796
$ws4.Range("C13").Value = 31.2
$ws4.Range("F23").Font.Color = 0

# --- Move the new sheet to the end of the tab strip and make it the active sheet/selection ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Move($null, $lastSheet)
$ws4b = $wb.Worksheets.Item("Tissue Weight Summary")
$ws4b.Activate()
$ws4b.Range("E6").Select()
